$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2201800.2
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 2752000.2
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 2752000.2
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -2752138.2
$ws.Range("H106").Value = 1621.0834
$ws.Range("J106").Value = 1947.75
$ws.Range("L106").Value = 1947.75
$ws.Range("N106").Value = -3209.75
$ws.Range("H112").Value = 10640697
$ws.Range("J112").Value = 11113386
$ws.Range("L112").Value = 33340158
$ws.Range("N112").Value = -33342374
$ws.Range("H132").Value = 3581.2
$ws.Range("I132").Value = 3526.2856
$ws.Range("K132").Value = 10578.8568
$ws.Range("M132").Value = -8048.856800000001
$ws.Range("H138").Value = 38463388
$ws.Range("I138").Value = 1267.6154
$ws.Range("K138").Value = 3802.8462
$ws.Range("M138").Value = 1337.1538
$ws.Range("H141").Value = 1146.2727
$ws.Range("I141").Value = 1019.9032
$ws.Range("J141").Value = 3105
$ws.Range("K141").Value = 3059.7096
$ws.Range("L141").Value = 9315
$ws.Range("M141").Value = 2120.2904
$ws.Range("N141").Value = -19675
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10641731
$ws.Range("J32").Value = 5913.8335
$ws.Range("L32").Value = 5913.8335
$ws.Range("N32").Value = -6487.8335
$ws.Range("H61").Value = 33337454
$ws.Range("I61").Value = 43481184
$ws.Range("J61").Value = 8054.4287
$ws.Range("K61").Value = 43481184
$ws.Range("L61").Value = 8054.4287
$ws.Range("M61").Value = -43480972
$ws.Range("N61").Value = -8478.4287
$ws.Range("H74").Value = 35757228
$ws.Range("I74").Value = 35757228
$ws.Range("K74").Value = 35757228
$ws.Range("M74").Value = -35756354
$ws.Range("H77").Value = 35757228
$ws.Range("I77").Value = 35757228
$ws.Range("K77").Value = 178786140
$ws.Range("M77").Value = -178781772
$ws.Range("H132").Value = 23313644
$ws.Range("I132").Value = 1618.697
$ws.Range("K132").Value = 4856.090999999999
$ws.Range("M132").Value = -2326.090999999999
$ws.Range("H136").Value = 33337454
$ws.Range("I136").Value = 43481184
$ws.Range("J136").Value = 8054.4287
$ws.Range("K136").Value = 130443552
$ws.Range("L136").Value = 24163.2861
$ws.Range("M136").Value = -130441002
$ws.Range("N136").Value = -29263.2861
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10018.462
$ws.Range("I105").Value = 10018.462
$ws.Range("K105").Value = 10018.462
$ws.Range("M105").Value = -8271.462
$ws.Range("H107").Value = 3248.3333
$ws.Range("I107").Value = 2766.4443
$ws.Range("J107").Value = 4694
$ws.Range("K107").Value = 2766.4443
$ws.Range("L107").Value = 4694
$ws.Range("M107").Value = -846.4443000000001
$ws.Range("N107").Value = -8534
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27783870
$ws.Range("I31").Value = 3698.5
$ws.Range("K31").Value = 3698.5
$ws.Range("M31").Value = -3403.5
$ws.Range("H34").Value = 27783870
$ws.Range("I34").Value = 3698.5
$ws.Range("K34").Value = 3698.5
$ws.Range("M34").Value = -3496.5
$ws.Range("H107").Value = 999.6667
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H134").Value = 4719.222
$ws.Range("I134").Value = 4519.5557
$ws.Range("K134").Value = 13558.6671
$ws.Range("M134").Value = -11023.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1412.1818
$ws.Range("I3").Value = 1053.4
$ws.Range("K3").Value = 3160.2
$ws.Range("M3").Value = -3048.2
$ws.Range("H113").Value = 1900.6
$ws.Range("J113").Value = 2248.3125
$ws.Range("L113").Value = 6744.9375
$ws.Range("N113").Value = -11084.9375
$ws.Range("H122").Value = 2016.1765
$ws.Range("J122").Value = 2085.0667
$ws.Range("L122").Value = 18765.6003
$ws.Range("N122").Value = -23665.6003
$ws.Range("H131").Value = 1637.0588
$ws.Range("I131").Value = 1429.8
$ws.Range("J131").Value = 1723.4166
$ws.Range("K131").Value = 4289.4
$ws.Range("L131").Value = 5170.2498
$ws.Range("M131").Value = 750.6000000000004
$ws.Range("N131").Value = -15250.2498
$ws.Range("H132").Value = 3032962.5
$ws.Range("J132").Value = 5131586.5
$ws.Range("L132").Value = 46184278.5
$ws.Range("N132").Value = -46189338.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1699.75
$ws.Range("I80").Value = 1266
$ws.Range("J80").Value = 3001
$ws.Range("K80").Value = 1266
$ws.Range("L80").Value = 3001
$ws.Range("M80").Value = -268
$ws.Range("N80").Value = -4997
$ws.Range("H83").Value = 1699.75
$ws.Range("I83").Value = 1266
$ws.Range("J83").Value = 3001
$ws.Range("K83").Value = 6330
$ws.Range("L83").Value = 15005
$ws.Range("M83").Value = -1338
$ws.Range("N83").Value = -24989
$ws.Range("H102").Value = 3300.6667
$ws.Range("J102").Value = 4997.3335
$ws.Range("L102").Value = 4997.3335
$ws.Range("N102").Value = -8241.333500000001
$ws.Range("H126").Value = 25087264
$ws.Range("I126").Value = 9186106
$ws.Range("J126").Value = 200000000
$ws.Range("K126").Value = 27558318
$ws.Range("L126").Value = 600000000
$ws.Range("M126").Value = -27555848
$ws.Range("N126").Value = -600004940
$ws.Range("H132").Value = 1494.7241
$ws.Range("I132").Value = 1433.4
$ws.Range("J132").Value = 1878
$ws.Range("K132").Value = 4300.200000000001
$ws.Range("L132").Value = 5634
$ws.Range("M132").Value = -1770.200000000001
$ws.Range("N132").Value = -10694
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3304.8948
$ws.Range("I16").Value = 3340.7646
$ws.Range("K16").Value = 3340.7646
$ws.Range("M16").Value = -3170.7646
$ws.Range("H22").Value = 1178.6097
$ws.Range("J22").Value = 1240.4166
$ws.Range("L22").Value = 1240.4166
$ws.Range("N22").Value = -1830.4166
$ws.Range("H27").Value = 1178.6097
$ws.Range("J27").Value = 1240.4166
$ws.Range("L27").Value = 1240.4166
$ws.Range("N27").Value = -1454.4166
$ws.Range("H46").Value = 1835.2941
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 6666.6665
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 6666.6665
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -7042.6665
$ws.Range("H132").Value = 125010584
$ws.Range("I132").Value = 6851.636
$ws.Range("J132").Value = 400018800
$ws.Range("K132").Value = 20554.908
$ws.Range("L132").Value = 1200056400
$ws.Range("M132").Value = -18024.908
$ws.Range("N132").Value = -1200061460
$ws.Range("H136").Value = 3988.25
$ws.Range("I136").Value = 3985.3125
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 11955.9375
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -9405.9375
$ws.Range("N136").Value = -17100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 876.2121
$ws.Range("I107").Value = 709.4583
$ws.Range("J107").Value = 1320.8889
$ws.Range("K107").Value = 2128.3749
$ws.Range("L107").Value = 3962.6667
$ws.Range("M107").Value = -208.3748999999998
$ws.Range("N107").Value = -7802.6667
$ws.Range("H132").Value = 1891.3256
$ws.Range("I132").Value = 1723.3334
$ws.Range("J132").Value = 2445.7
$ws.Range("K132").Value = 5170.0002
$ws.Range("L132").Value = 7337.099999999999
$ws.Range("M132").Value = -2640.0002
$ws.Range("N132").Value = -12397.1
$ws.Range("H136").Value = 9699.714
$ws.Range("I136").Value = 2999.5
$ws.Range("J136").Value = 12379.8
$ws.Range("K136").Value = 8998.5
$ws.Range("L136").Value = 37139.39999999999
$ws.Range("M136").Value = -6448.5
$ws.Range("N136").Value = -42239.39999999999
